# "Added Armor of Vitality!"
#
# Insert a new card row into the Fortitude sheet: "Armor of Vitality",
# Type=Skill, Cost=1, Rarity=UNCOMMON. This becomes row 2 (right after
# "Ablative Skin"), pushing every subsequent card down by one row.
# Downstream COUNTIF/SUM summary formulas on Sheet1 recalculate on their
# own once the new row is in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fortitude")

[void]$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "Armor of Vitality"
$ws.Cells.Item(2, 2).Value = "Skill"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = "UNCOMMON"

# Restore the selection / active-sheet bookkeeping recorded in the diff:
# Sheet1's selection moved to D16, Celerity is no longer the active tab,
# and Fortitude (now with cell A1 selected) became the active tab instead.
$ws1 = $wb.Worksheets.Item("Sheet1")
[void]$ws1.Activate()
[void]$ws1.Range("D16").Select()

$wsCelerity = $wb.Worksheets.Item("Celerity")
[void]$wsCelerity.Activate()
[void]$wsCelerity.Range("A15").Select()

[void]$ws.Activate()
[void]$ws.Range("A1").Select()
